$d = $word.ActiveDocument

# --- Change 1: add first-line indent to the title paragraph ---
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.FirstLineIndent = 36

# --- Change 2/3: move the "_GoBack" bookmark so it wraps the text that
#     starts at "At a minimum" and ends at "...deem useful." in the
#     Metadata section. Word only allows one bookmark per name, so
#     re-adding "_GoBack" here automatically removes it from its old
#     location at the end of the document. ---
$startRng = $d.Content
[void]$startRng.Find.Execute("At a minimum, you need the dimensions")
$startPos = $startRng.Start

$endRng = $d.Content
[void]$endRng.Find.Execute("and time, but you can include any information you deem useful.")
$endPos = $endRng.End

$bmRange = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "done"
